$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the "5840793 - Sérgio Schneider" text in columns B/C (old
# row 13, with nothing in column A) is removed; everything below it shifts up by one.
$ws.Rows("13:13").Delete()

# After the shift, several cells (duplicated across columns B and C) need their text
# content replaced with the new values.
$ws.Range("B10").Value = "5840793 - Sérgio Schneider"
$ws.Range("C10").Value = "5840793 - Sérgio Schneider"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2012" looks like a date, and a plain .Value assignment would get converted
# to a date serial number. Copy the existing text value from A8's row (which already
# holds this exact string as text) so the new cell keeps its text type and original
# style.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("B18").Value = "5840793 - Sérgio Schneider"
$ws.Range("C18").Value = "5840793 - Sérgio Schneider"

$ws.Range("B19").Value = "Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios."
$ws.Range("C19").Value = "Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios."

$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
